# Commit: "fixed feedback loop + added connection to player's moves"
#
# Content change: the PrioritizationAlgorithm parameter (column K) on the
# "ScenarioParameters" sheet was changed from {5, 4, 4} to {2, 2, 2} for the
# three scenario rows - this is the "fixed feedback loop" part of the commit.

$wb = $excel.ActiveWorkbook

$wsInfo   = $wb.Worksheets.Item("ScenarioInfo")
$wsParams = $wb.Worksheets.Item("ScenarioParameters")

# --- Data fix: PrioritizationAlgorithm (K2:K4) 5,4,4 -> 2,2,2 ---
$wsParams.Range("K2").Value = 2
$wsParams.Range("K3").Value = 2
$wsParams.Range("K4").Value = 2

# --- New column (AutoIntensificationKM, column L) gained an explicit width ---
$wsParams.Columns.Item(12).ColumnWidth = 17.86

# --- Window / selection state: user ended up with ScenarioParameters active,
#     having last looked at cell C4, after previously looking at C16 on
#     ScenarioInfo (which is no longer the selected tab). SpecsData is left
#     untouched. ---
$wsInfo.Activate()
$wsInfo.Range("C16").Select()

$wsParams.Activate()
$wsParams.Range("C4").Select()
